$wb = $excel.ActiveWorkbook

# Rename the "skillId" column header to "levelPackId" on the ActorLevelPackTable sheet
$ws = $wb.Worksheets.Item("ActorLevelPackTable")
$ws.Range("B1").Value = "levelPackId|String"

# Widen column B so the longer header text ("levelPackId|String") is no
# longer clipped (mirrors the author widening the column after retyping
# the header).
$ws.Columns.Item(2).ColumnWidth = 16.857142857142858

# The author was working in ActorLevelPackTable when saving, so it ends
# up as the active sheet/tab.
$ws.Activate()
